# Applies the commit "Popravke testova i dopuna podacima data" to Sheet1.
#
# Logical changes (reconstructed from the OOXML diff):
#   1. Header K1: "current_account_3_bban" -> "current_account_3_iban"
#   2. I2 (current_account_2_iban): "." -> "RS35 2059 0310 0441 9532 81"
#   3. K2 (current_account_3_iban, was _bban): "205-9031004419532-81" -> "."
#      (the now-unreferenced shared strings "current_account_3_bban" and
#      "205-9031004419532-81" simply disappear from sharedStrings.xml once
#      nothing points at them any more - the engine regenerates that table
#      on save, so no direct action is needed for that part.)
#   4. Sheet view: active selection moves from M1 to M22.
#   5. Column K (11) gets its own explicit width (was merged with column J
#      as one <col min="10" max="11">, now split so K is wider).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2 & 3: cell content updates -----------------------------------
$ws.Range("K1").Value = "current_account_3_iban"
$ws.Range("I2").Value = "RS35 2059 0310 0441 9532 81"
$ws.Range("K2").Value = "."

# --- 4: column K width, split away from column J -----------------------
$ws.Columns.Item(11).ColumnWidth = 27

# --- 5: move the active selection to M22 --------------------------------
$ws.Range("M22").Select()
